$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1. Insert a new column at D (keeps the existing "measure_name" column's
#    formatting in place at C, and pushes D..H -> E..I). We then copy the
#    old column C (measure_name) values into the new D column, and finally
#    overwrite column C with the new "measure_etl_name" values. Doing the
#    copy BEFORE the overwrite keeps the shared string referenced so it
#    isn't garbage collected out of the shared-strings table.
# ---------------------------------------------------------------------------
$ws.Columns.Item(4).Insert()
$ws.Range("C1:C30").Copy()
$ws.Range("D1").PasteSpecial(-4163)
$ws.Application.CutCopyMode = 0

# New column D width matches column C's width value (without the bestFit flag)
$ws.Columns.Item(4).ColumnWidth = 53

# ---------------------------------------------------------------------------
# 2. Populate column C (measure_etl_name) header + values.
# ---------------------------------------------------------------------------
$ws.Range("C1").Value = "measure_etl_name"

$ws.Range("C2").Value = "All-Cause ED Visits"
$ws.Range("C3").Value = "Acute Hospital Utilization"
$ws.Range("C4").Value = "Follow-up ED visit for Alcohol/Drug Abuse"
$ws.Range("C5").Value = "Follow-up ED visit for Alcohol/Drug Abuse"
$ws.Range("C6").Value = "Follow-up ED visit for Mental Illness"
$ws.Range("C7").Value = "Follow-up ED visit for Mental Illness"
$ws.Range("C8").Value = "Follow-up Hospitalization for Mental Illness"
$ws.Range("C9").Value = "Follow-up Hospitalization for Mental Illness"
$ws.Range("C10").Value = "Mental Health Treatment Penetration"
$ws.Range("C11").Value = "SUD Treatment Penetration"
$ws.Range("C12").Value = "SUD Treatment Penetration (Opioid)"
$ws.Range("C13").Value = "Plan All-Cause Readmissions (30 days)"
$ws.Range("C14").Value = "Child and Adolescent Access to Primary Care"
$ws.Range("C15").Value = "Diabetes Care: Eye Exam"
$ws.Range("C16").Value = "Diabetes Care: A1c Testing"
$ws.Range("C17").Value = "Diabetes Care: Kidney Screening"
$ws.Range("C18").Value = "Medication Management for Asthma: Compliance 50%"
$ws.Range("C19").Value = "Medication Management for Asthma: Compliance 75%"
$ws.Range("C20").Value = "Asthma Medication Ratio"
$ws.Range("C21").Value = "Asthma Medication Ratio (1-year requirement)"
$ws.Range("C22").Value = "Percent Homeless"
$ws.Range("C23").Value = "Antidepressant Medication Management"
$ws.Range("C24").Value = "High-dose Chronic Opioid Therapy"
$ws.Range("C25").Value = "Concurrent Opioids and Sedatives Prescriptions"
$ws.Range("C26").Value = "Statin Therapy for Heart Disease"
$ws.Range("C27").Value = "SUD Treatment Initiation"
$ws.Range("C28").Value = "SUD Treatment Initiation"
$ws.Range("C29").Value = "SUD Treatment Initiation (No Modifiers)"
$ws.Range("C30").Value = "SUD Treatment Initiation (No Modifiers)"

# ---------------------------------------------------------------------------
# 3. Append the new "MH Treatment Penetration by Diagnosis" rows (31-37).
# ---------------------------------------------------------------------------
$newRows = @(
  @(30, "TPM_ADHD",       "MH Treatment Penetration: ADHD"),
  @(31, "TPM_Adjustment", "MH Treatment Penetration: Adjustment"),
  @(32, "TPM_Anxiety",    "MH Treatment Penetration: Anxiety"),
  @(33, "TPM_Depression", "MH Treatment Penetration: Depression"),
  @(34, "TPM_Impulse",    "MH Treatment Penetration: Disrup/Impulse/Conduct"),
  @(35, "TPM_Bipolar",    "MH Treatment Penetration: Mania/Bipolar"),
  @(36, "TPM_Psychotic",  "MH Treatment Penetration: Psychotic")
)

$r = 31
foreach ($row in $newRows) {
    $ws.Cells.Item($r, 1).Value = $row[0]
    $ws.Cells.Item($r, 2).Value = $row[1]
    $ws.Cells.Item($r, 3).Value = "MH Treatment Penetration by Diagnosis"
    $ws.Cells.Item($r, 4).Value = $row[2]
    $ws.Cells.Item($r, 5).Value = "age_grp_5"
    $ws.Cells.Item($r, 6).Value = "Age 6-17, Age 18-64, Age 65+"
    $ws.Cells.Item($r, 8).Value = "Proportion of members"
    $r = $r + 1
}

# ---------------------------------------------------------------------------
# 4. Formulas: I2 is a standalone formula, I3:I37 is a shared formula group.
# ---------------------------------------------------------------------------
$f2 = '=CONCATENATE(",(",A2,", ''",B2,"'', ''",C2,"'', ''",D2,"''",", ''",E2,"'', ''",F2,"'', ''",G2,"'', ''",H2,"'')")'
$ws.Range("I2").Formula = $f2

$f = '=CONCATENATE(",(",A3,", ''",B3,"'', ''",C3,"'', ''",D3,"''",", ''",E3,"'', ''",F3,"'', ''",G3,"'', ''",H3,"'')")'
$ws.Range("I3:I37").Formula = $f
